# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" between "2021-Q3" and "总计", populated
#    with the per-fund holdings detail for the new quarter.
# 2) Update the "总计" (totals) sheet with a new leading row summarising the
#    2022-Q1 quarter, pushing the existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet, positioned right after "2021-Q3"
# ---------------------------------------------------------------------------
# NOTE: sheet handles returned by Worksheets.Item(...) are position-bound, so
# any reference taken *before* the insertion (e.g. to "总计", which sits
# after the insertion point) would silently start pointing at the newly
# inserted sheet once positions shift. Fetch "总计" again after inserting.
$q3 = $wb.Worksheets.Item("2021-Q3")

$newWs = $wb.Worksheets.Add($null, $q3)
$newWs.Name = "2022-Q1"

# Match the page margins used by the sibling quarterly sheets (0.75"/1"/0.5").
$newWs.PageSetup.LeftMargin = 54
$newWs.PageSetup.RightMargin = 54
$newWs.PageSetup.TopMargin = 72
$newWs.PageSetup.BottomMargin = 72
$newWs.PageSetup.HeaderMargin = 36
$newWs.PageSetup.FooterMargin = 36

$total = $wb.Worksheets.Item("总计")

# Match the look & feel of the other quarterly sheets: copy the bold/centered
# header style (columns B:H on row 1) and the index-column style (column A)
# from the "2021-Q3" sheet, which already carries the formatting we want.
$q3.Range("B1:H1").Copy()
$newWs.Range("B1:H1").PasteSpecial(-4122)

$q3.Range("A2").Copy()
$newWs.Range("A2:A4").PasteSpecial(-4122)

# Headers
$newWs.Range("B1").Value = "基金代码"
$newWs.Range("C1").Value = "基金名称"
$newWs.Range("D1").Value = "基金规模"
$newWs.Range("E1").Value = "股票总仓位"
$newWs.Range("F1").Value = "仓位占比"
$newWs.Range("G1").Value = "持有市值(亿元)"
$newWs.Range("H1").Value = "仓位排名"

# Numeric-looking identifiers / metrics are stored as text in this workbook
# (e.g. "519613", "5.70") -- force text formatting before assigning so they
# round-trip as strings rather than being coerced to numbers.
$newWs.Range("B2:B4").NumberFormat = "@"
$newWs.Range("D2:G4").NumberFormat = "@"

# Row 2 - 银河君尚灵活配置混合A
$newWs.Range("A2").Value = 0
$newWs.Range("B2").Value = "519613"
$newWs.Range("C2").Value = "银河君尚灵活配置混合A"
$newWs.Range("D2").Value = "5.70"
$newWs.Range("E2").Value = "29.99"
$newWs.Range("F2").Value = "0.47"
$newWs.Range("G2").Value = "0.0268"
$newWs.Range("H2").Value = 9

# Row 3 - 银河君尚灵活配置混合I
$newWs.Range("A3").Value = 1
$newWs.Range("B3").Value = "519615"
$newWs.Range("C3").Value = "银河君尚灵活配置混合I"
$newWs.Range("D3").Value = "5.70"
$newWs.Range("E3").Value = "29.99"
$newWs.Range("F3").Value = "0.47"
$newWs.Range("G3").Value = "0.0268"
$newWs.Range("H3").Value = 9

# Row 4 - 银河君尚灵活配置混合C
$newWs.Range("A4").Value = 2
$newWs.Range("B4").Value = "519614"
$newWs.Range("C4").Value = "银河君尚灵活配置混合C"
$newWs.Range("D4").Value = "0.23"
$newWs.Range("E4").Value = "29.99"
$newWs.Range("F4").Value = "0.47"
$newWs.Range("G4").Value = "0.0011"
$newWs.Range("H4").Value = 9

# ---------------------------------------------------------------------------
# Step 2: push the "总计" rows down and insert the new 2022-Q1 summary row
# ---------------------------------------------------------------------------

# Remember the current (pre-edit) contents of rows 2 and 3 before they are
# overwritten, so they can be shifted down into rows 3 and 4.
$oldRow2Date  = $total.Range("B2").Value()
$oldRow2Count = $total.Range("C2").Value()
$oldRow2Value = $total.Range("D2").Value()

$oldRow3Date  = $total.Range("B3").Value()
$oldRow3Count = $total.Range("C3").Value()
$oldRow3Value = $total.Range("D3").Value()

# Extend the bold index-column style (currently only on A2) down to A4.
$total.Range("A2").Copy()
$total.Range("A4").PasteSpecial(-4122)

# Row 4 (was row 3): 2021-Q2
$total.Range("A4").Value = 2
$total.Range("B4").Value = $oldRow3Date
$total.Range("C4").Value = $oldRow3Count
$total.Range("D4").Value = $oldRow3Value

# Row 3 (was row 2): 2021-Q3
$total.Range("A3").Value = 1
$total.Range("B3").Value = $oldRow2Date
$total.Range("C3").Value = $oldRow2Count
$total.Range("D3").Value = $oldRow2Value

# Row 2 (new): 2022-Q1
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.05

# Restore the originally active sheet/selection (adding a sheet shifts focus
# to it by default) so the workbook-level view state is left untouched.
$wb.Worksheets.Item("2021-Q2").Activate()
[void]$wb.Worksheets.Item("2021-Q2").Range("A1").Select()

Write-Host "2022-Q1 sheet added and totals updated"
